$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that resolved as "Fallo" (loss): resultado = "Fallo", profit = -1
$falloRows = @(292, 295, 297, 298, 299, 300, 310, 315)
foreach ($r in $falloRows) {
    $ws.Range("G$r").Value = "Fallo"
    $ws.Range("H$r").Value = -1
}

# Rows that resolved as "Acierto" (win): resultado = "Acierto", profit = cuota - 1
$aciertoRows = @{
    312 = 0.91
    313 = 1.75
    314 = 1.62
    317 = 0.8
    322 = 1.5
    324 = 2.25
}
foreach ($r in $aciertoRows.Keys) {
    $ws.Range("G$r").Value = "Acierto"
    $ws.Range("H$r").Value = $aciertoRows[$r]
}

# event_id for row 325 was stored as text; convert it to a true number
$ws.Range("A325").Value = 14466743
